$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 20 new rows before old row 21 to make room for new
# "komorbiditet ICD10 symbol" row and 19 "CCI" rows, shifting
# everything from the old row 21 onward down to row 41.
$ws.Rows("21:40").Insert()

# Update row 20 (komorbiditet / items) value to the new, shorter dict
$ws.Cells.Item(20, 4).Value2 = '{''Sykdomskategori'': ''Fordøyelsessystemet, sykdommer (Kap.XI; K00-K93)'', ''Sykdomstilstand'': ''Psykiske lidelser og atferdsforstyrrelser (Kap.V; F00-F99)''}'

# Populate the 20 newly-inserted rows (21-40)
$ws.Cells.Item(21, 1).Value2 = 'komorbiditet'
$ws.Cells.Item(21, 2).Value2 = 'ICD10 symbol'
$ws.Cells.Item(21, 3).Value2 = 1
$ws.Cells.Item(21, 4).Value2 = 'diam tincidunt erat marti....'

$ws.Cells.Item(22, 1).Value2 = 'CCI'
$ws.Cells.Item(22, 2).Value2 = 'items'
$ws.Cells.Item(22, 3).Value2 = 1
$ws.Cells.Item(22, 4).Value2 = '{''Kategori'': ''Aldersgruppe'', ''Verdi'': ''1'', ''Symbol'': ''50-59 år''}'

$ws.Cells.Item(23, 1).Value2 = 'CCI'
$ws.Cells.Item(23, 2).Value2 = 'items'
$ws.Cells.Item(23, 3).Value2 = 2
$ws.Cells.Item(23, 4).Value2 = '{''Kategori'': ''Myokardinfarkt'', ''Verdi'': ''0'', ''Symbol'': ''Nei''}'

$ws.Cells.Item(24, 1).Value2 = 'CCI'
$ws.Cells.Item(24, 2).Value2 = 'items'
$ws.Cells.Item(24, 3).Value2 = 3
$ws.Cells.Item(24, 4).Value2 = '{''Kategori'': ''Kronisk hjertesvikt'', ''Verdi'': ''0'', ''Symbol'': ''Nei''}'

$ws.Cells.Item(25, 1).Value2 = 'CCI'
$ws.Cells.Item(25, 2).Value2 = 'items'
$ws.Cells.Item(25, 3).Value2 = 4
$ws.Cells.Item(25, 4).Value2 = '{''Kategori'': ''Perifere vaskulære sykdommer'', ''Verdi'': ''0'', ''Symbol'': ''Nei''}'

$ws.Cells.Item(26, 1).Value2 = 'CCI'
$ws.Cells.Item(26, 2).Value2 = 'items'
$ws.Cells.Item(26, 3).Value2 = 5
$ws.Cells.Item(26, 4).Value2 = '{''Kategori'': ''Cerebrovaskulær hendelse'', ''Verdi'': ''0'', ''Symbol'': ''Nei''}'

$ws.Cells.Item(27, 1).Value2 = 'CCI'
$ws.Cells.Item(27, 2).Value2 = 'items'
$ws.Cells.Item(27, 3).Value2 = 6
$ws.Cells.Item(27, 4).Value2 = '{''Kategori'': ''Demens'', ''Verdi'': ''0'', ''Symbol'': ''Nei''}'

$ws.Cells.Item(28, 1).Value2 = 'CCI'
$ws.Cells.Item(28, 2).Value2 = 'items'
$ws.Cells.Item(28, 3).Value2 = 7
$ws.Cells.Item(28, 4).Value2 = '{''Kategori'': ''Kronisk lungesykdom'', ''Verdi'': ''0'', ''Symbol'': ''Nei''}'

$ws.Cells.Item(29, 1).Value2 = 'CCI'
$ws.Cells.Item(29, 2).Value2 = 'items'
$ws.Cells.Item(29, 3).Value2 = 8
$ws.Cells.Item(29, 4).Value2 = '{''Kategori'': ''Ulcussykdom'', ''Verdi'': ''0'', ''Symbol'': ''Nei''}'

$ws.Cells.Item(30, 1).Value2 = 'CCI'
$ws.Cells.Item(30, 2).Value2 = 'items'
$ws.Cells.Item(30, 3).Value2 = 9
$ws.Cells.Item(30, 4).Value2 = '{''Kategori'': ''Leversykdom'', ''Verdi'': ''0'', ''Symbol'': ''Ingen''}'

$ws.Cells.Item(31, 1).Value2 = 'CCI'
$ws.Cells.Item(31, 2).Value2 = 'items'
$ws.Cells.Item(31, 3).Value2 = 10
$ws.Cells.Item(31, 4).Value2 = '{''Kategori'': ''Bindevevssykdom'', ''Verdi'': ''0'', ''Symbol'': ''Nei''}'

$ws.Cells.Item(32, 1).Value2 = 'CCI'
$ws.Cells.Item(32, 2).Value2 = 'items'
$ws.Cells.Item(32, 3).Value2 = 11
$ws.Cells.Item(32, 4).Value2 = '{''Kategori'': ''Diabetes'', ''Verdi'': ''0'', ''Symbol'': ''Ingen eller diettkontrollert''}'

$ws.Cells.Item(33, 1).Value2 = 'CCI'
$ws.Cells.Item(33, 2).Value2 = 'items'
$ws.Cells.Item(33, 3).Value2 = 12
$ws.Cells.Item(33, 4).Value2 = '{''Kategori'': ''Hemiplegi'', ''Verdi'': ''0'', ''Symbol'': ''Nei''}'

$ws.Cells.Item(34, 1).Value2 = 'CCI'
$ws.Cells.Item(34, 2).Value2 = 'items'
$ws.Cells.Item(34, 3).Value2 = 13
$ws.Cells.Item(34, 4).Value2 = '{''Kategori'': ''Moderat til alvorlig nyresykdom'', ''Verdi'': ''0'', ''Symbol'': ''Nei''}'

$ws.Cells.Item(35, 1).Value2 = 'CCI'
$ws.Cells.Item(35, 2).Value2 = 'items'
$ws.Cells.Item(35, 3).Value2 = 14
$ws.Cells.Item(35, 4).Value2 = '{''Kategori'': ''Solid svulst'', ''Verdi'': ''0'', ''Symbol'': ''Ingen''}'

$ws.Cells.Item(36, 1).Value2 = 'CCI'
$ws.Cells.Item(36, 2).Value2 = 'items'
$ws.Cells.Item(36, 3).Value2 = 15
$ws.Cells.Item(36, 4).Value2 = '{''Kategori'': ''Leukemi'', ''Verdi'': ''0'', ''Symbol'': ''Nei''}'

$ws.Cells.Item(37, 1).Value2 = 'CCI'
$ws.Cells.Item(37, 2).Value2 = 'items'
$ws.Cells.Item(37, 3).Value2 = 16
$ws.Cells.Item(37, 4).Value2 = '{''Kategori'': ''Lymfom'', ''Verdi'': ''0'', ''Symbol'': ''Nei''}'

$ws.Cells.Item(38, 1).Value2 = 'CCI'
$ws.Cells.Item(38, 2).Value2 = 'items'
$ws.Cells.Item(38, 3).Value2 = 17
$ws.Cells.Item(38, 4).Value2 = '{''Kategori'': ''AIDS'', ''Verdi'': ''0'', ''Symbol'': ''Nei''}'

$ws.Cells.Item(39, 1).Value2 = 'CCI'
$ws.Cells.Item(39, 2).Value2 = 'items'
$ws.Cells.Item(39, 3).Value2 = 18
$ws.Cells.Item(39, 4).Value2 = '{''Kategori'': ''CCI totalskår'', ''Verdi'': ''82''}'

$ws.Cells.Item(40, 1).Value2 = 'CCI'
$ws.Cells.Item(40, 2).Value2 = 'items'
$ws.Cells.Item(40, 3).Value2 = 19
$ws.Cells.Item(40, 4).Value2 = '{''Kategori'': ''Estimert overlevelse etter 10 år'', ''Verdi'': ''82'', ''Enhet'': ''%''}'

